# deep sea double count fix
# Update recalculated values on three related sheets after fixing a
# double-counting bug in the "deep sea" landings aggregation.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Status by Landings (Area)" ---
$wsArea = $wb.Worksheets.Item("Status by Landings (Area)")
$wsArea.Range("C3").Value  = 0.9823068853471902
$wsArea.Range("C4").Value  = 0.2254840312804062
$wsArea.Range("C5").Value  = 1.197280435915919
$wsArea.Range("C6").Value  = 0.2254840312804062
$wsArea.Range("C7").Value  = 15.10956701022704
$wsArea.Range("C8").Value  = 69.04212945962217
$wsArea.Range("C9").Value  = 15.84830353015079
$wsArea.Range("C10").Value = 84.15169646984921
$wsArea.Range("C11").Value = 15.84830353015079

# --- Sheet: "Status by Landings (Tier)" ---
$wsTier = $wb.Worksheets.Item("Status by Landings (Tier)")

# Row 4 ("21")
$wsTier.Range("C4").Value = 0.9823068853471902
$wsTier.Range("D4").Value = 0.2254840312804062
$wsTier.Range("E4").Value = 1.197280435915919
$wsTier.Range("F4").Value = 0.2254840312804062
$wsTier.Range("G4").Value = 15.10956701022704
$wsTier.Range("H4").Value = 69.04212945962217
$wsTier.Range("I4").Value = 15.84830353015079
$wsTier.Range("J4").Value = 84.15169646984921
$wsTier.Range("K4").Value = 15.84830353015079

# Row 5 ("Global")
$wsTier.Range("C5").Value = 0.9823068853471902
$wsTier.Range("D5").Value = 0.2254840312804062
$wsTier.Range("E5").Value = 1.197280435915919
$wsTier.Range("F5").Value = 0.2254840312804062
$wsTier.Range("G5").Value = 15.10956701022704
$wsTier.Range("H5").Value = 69.04212945962217
$wsTier.Range("I5").Value = 15.84830353015079
$wsTier.Range("J5").Value = 84.15169646984921
$wsTier.Range("K5").Value = 15.84830353015079

# --- Sheet: "Comparison by Landings" ---
$wsComp = $wb.Worksheets.Item("Comparison by Landings")
$wsComp.Range("C2").Value = 92.31243281328341
$wsComp.Range("C3").Value = 15.10956701022704
$wsComp.Range("C4").Value = 69.04212945962217
$wsComp.Range("C5").Value = 15.84830353015079
$wsComp.Range("C6").Value = 84.15169646984921
$wsComp.Range("C7").Value = 15.84830353015079
